$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("B3").Value = 0.35
$ws.Range("C3").Value = 0.15

$ws.Range("B4").Value = 0.09
$ws.Range("C4").Value = 0.04

$ws.Range("B5").Value = 0.05
$ws.Range("C5").Value = 0.04

$ws.Range("B7").Value = 0.19

$ws.Range("B8").Value = 0.08

$ws.Range("B11").Value = 0.35
